$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update the Date value (row 8, column B) ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2026-02-05T08:09:31+00:00"

# --- Sheet "Elements": rename "topographique" -> "precisionTopographique" ---
$wsElem = $wb.Worksheets.Item("Elements")

# Row 13: fr-lm-technique-imagerie.lateralite.topographique
$wsElem.Range("A13").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique"
$wsElem.Range("B13").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique"
$wsElem.Range("AF13").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique"

# Row 14: fr-lm-technique-imagerie.lateralite.topographique.id
$wsElem.Range("A14").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.id"
$wsElem.Range("B14").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.id"

# Row 15: fr-lm-technique-imagerie.lateralite.topographique.extension
$wsElem.Range("A15").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.extension"
$wsElem.Range("B15").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.extension"

# Row 16: fr-lm-technique-imagerie.lateralite.topographique.coding
$wsElem.Range("A16").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.coding"
$wsElem.Range("B16").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.coding"

# Row 17: fr-lm-technique-imagerie.lateralite.topographique.text
$wsElem.Range("A17").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.text"
$wsElem.Range("B17").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.text"

# Row 18: fr-lm-technique-imagerie.lateralite.topographique.nom
$wsElem.Range("A18").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.nom"
$wsElem.Range("B18").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.nom"
$wsElem.Range("AF18").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.nom"

# Row 19: fr-lm-technique-imagerie.lateralite.topographique.valeur
$wsElem.Range("A19").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.valeur"
$wsElem.Range("B19").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.valeur"
$wsElem.Range("AF19").Value = "fr-lm-technique-imagerie.lateralite.precisionTopographique.valeur"

# --- Column widths on "Elements" sheet, widened to fit the longer text ---
$wsElem.Columns.Item(1).ColumnWidth = 55.7890625
$wsElem.Columns.Item(2).ColumnWidth = 55.7890625
$wsElem.Columns.Item(32).ColumnWidth = 53.0390625
